# Commit message: "renamed repo, fixed output folder path"
# The underlying data pipeline changed which per-genome FASTA labels were
# included for two species folders, which removed some now-duplicate /
# stale rows from each prediction sheet. Replicate that by deleting the
# affected rows from each worksheet (remaining rows shift up automatically).

$wb = $excel.ActiveWorkbook

# Sheet 1: "s__UBA1685 sp002320595-b-p" -- remove old rows 9-12
# (label_UMGS106_11.fasta, label_UMGS106_18.fasta, label_UMGS106_2.fasta,
#  label_UMGS106_3.fasta), dimension shrinks from A1:E32 to A1:E28.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A9:A12").EntireRow.Delete()

# Sheet 2: "s__UBA1685 sp900546845-b-p" -- remove old rows 11-22
# (label_UMGS922_19.fasta, _25, _26, _3, _34, _36, _41, _42, _51, _54,
#  _55, _57.fasta), dimension shrinks from A1:E69 to A1:E57.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A11:A22").EntireRow.Delete()
